$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extra Wishlist")

$ws.Range("C20").Value = "https://www.apple.com/ch-de/shop/buy-watch/apple-watch/42mm-gps-ros%C3%A9gold-aluminium-mandarine-sportarmband-sm"
$ws.Range("A20").Value = "Apple Watch"
$ws.Range("B20").Value = "https://store.storeimages.cdn-apple.com/1/as-images.apple.com/is/MDRY4ref_VW_34FR+watch-case-42-aluminum-rosegold-nc-s10_VW_34FR+watch-face-42-aluminum-rosegold-s10_VW_34FR?wid=5120&hei=3280&bgc=fafafa&trim=1&fmt=p-jpg&qlt=80&.v=aWR3eURSMEloeW8rRlJjeisrTUVuSjNqbENGcEFFTlJVaXJwL2VzdGxEM3hoaFEySTlrZSsrbEhqczIxajV3MzBnRnJNVU50TjZtbUN3ZjYvVksrL1FtTHNhRmRKQkE2OXJHRVEvRTFWdXRGd1ZOdUM2RkNJOUpIeklzU1JHbVNWNEM0ZHZzMzMzajJiOGZuamRiTE9rUVN3R3VxZWhYYXgwOHljYmZFMXBocmMyRTN3NCt6QkoxaUdRb0FBay9VYktGTHdENW9lYUFnak5pcy9ReEdDV1huc2RkMVJhUksxbWsxK1ZQMHVvWVdkR0c1QU9HSFhRTXZyMjRWUXMzUA"
$ws.Range("D20").Value = 400

$ws.Range("D21").Select()
